$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$medium = -4138
$thin = 2
$xlLeft = 7
$xlTop = 8
$xlBottom = 9
$xlRight = 10
$xlCenter = -4108

function Set-Edge {
    param($rng, $edgeIdx, $weight)
    $rng.Borders.Item($edgeIdx).LineStyle = 1
    $rng.Borders.Item($edgeIdx).Weight = $weight
}

function Set-ThinBox {
    param($rng)
    Set-Edge $rng $xlLeft $thin
    Set-Edge $rng $xlRight $thin
    Set-Edge $rng $xlTop $thin
    Set-Edge $rng $xlBottom $thin
}

# ---------------------------------------------------------------
# New small "inhalation model" lookup table, rows 85 (spacer) - 87
# ---------------------------------------------------------------

# Spacer row before the table
$ws.Rows.Item(85).RowHeight = 15.75

# ---- Row 87: data row (filled first so shared strings land in the
#      same order as the authored workbook: A, B, C, I, then the
#      header row A..M) ----

$cA87 = $ws.Range("A87")
$cA87.Value = "'104-55-2"
Set-ThinBox $cA87

$cB87 = $ws.Range("B87")
$cB87.Value = "3-phenyl-2-propenal"
Set-ThinBox $cB87

$cC87 = $ws.Range("C87")
$cC87.Value = "F"
Set-ThinBox $cC87
$cC87.HorizontalAlignment = $xlCenter

$cD87 = $ws.Range("D87")
Set-ThinBox $cD87
$cD87.HorizontalAlignment = $xlCenter

$cE87 = $ws.Range("E87")
Set-ThinBox $cE87
$cE87.HorizontalAlignment = $xlCenter

$cF87 = $ws.Range("F87")
Set-ThinBox $cF87
$cF87.HorizontalAlignment = $xlCenter

$cG87 = $ws.Range("G87")
Set-ThinBox $cG87
$cG87.HorizontalAlignment = $xlCenter

$cH87 = $ws.Range("H87")
Set-ThinBox $cH87
$cH87.HorizontalAlignment = $xlCenter

$cI87 = $ws.Range("I87")
$cI87.Value = "Sa"
Set-ThinBox $cI87
$cI87.HorizontalAlignment = $xlCenter

$cJ87 = $ws.Range("J87")
$cJ87.Value = 873.10334785411487
Set-ThinBox $cJ87
$cJ87.NumberFormat = "0.00E+00"

$cK87 = $ws.Range("K87")
$cK87.Value = 274.83920626249699
Set-ThinBox $cK87
$cK87.NumberFormat = "0.00E+00"

$cL87 = $ws.Range("L87")
$cL87.Value = 80886.100486578129
Set-ThinBox $cL87
$cL87.NumberFormat = "0.00E+00"

$cM87 = $ws.Range("M87")
$cM87.Value = 294.30335499268068
$cM87.NumberFormat = "0.00E+00"

# ---- Row 86: header row ----

$ws.Rows.Item(86).RowHeight = 120

$cA86 = $ws.Range("A86")
$cA86.Value = "CAS-number"
$cA86.Font.Bold = $true
Set-Edge $cA86 $xlLeft $medium
Set-Edge $cA86 $xlTop $medium
$cA86.HorizontalAlignment = $xlCenter
$cA86.WrapText = $true

$cB86 = $ws.Range("B86")
$cB86.Value = "Compound name"
$cB86.Font.Bold = $true
Set-Edge $cB86 $xlLeft $medium
Set-Edge $cB86 $xlRight $medium
Set-Edge $cB86 $xlTop $medium
$cB86.HorizontalAlignment = $xlCenter
$cB86.WrapText = $true

$cC86 = $ws.Range("C86")
$cC86.Value = "Faeces"
$cC86.Font.Bold = $true
Set-Edge $cC86 $xlTop $medium
$cC86.HorizontalAlignment = $xlCenter
$cC86.WrapText = $true

$cD86 = $ws.Range("D86")
$cD86.Value = "Urine"
$cD86.Font.Bold = $true
Set-Edge $cD86 $xlLeft $medium
Set-Edge $cD86 $xlRight $medium
Set-Edge $cD86 $xlTop $medium
$cD86.HorizontalAlignment = $xlCenter
$cD86.WrapText = $true

$cE86 = $ws.Range("E86")
$cE86.Value = "Breath"
$cE86.Font.Bold = $true
Set-Edge $cE86 $xlTop $medium
$cE86.HorizontalAlignment = $xlCenter
$cE86.WrapText = $true

$cF86 = $ws.Range("F86")
$cF86.Value = "Skin"
$cF86.Font.Bold = $true
Set-Edge $cF86 $xlLeft $medium
Set-Edge $cF86 $xlRight $medium
Set-Edge $cF86 $xlTop $medium
$cF86.HorizontalAlignment = $xlCenter
$cF86.WrapText = $true

$cG86 = $ws.Range("G86")
$cG86.Value = "Milk"
$cG86.Font.Bold = $true
Set-Edge $cG86 $xlTop $medium
$cG86.HorizontalAlignment = $xlCenter
$cG86.WrapText = $true

$cH86 = $ws.Range("H86")
$cH86.Value = "Blood"
$cH86.Font.Bold = $true
Set-Edge $cH86 $xlLeft $medium
Set-Edge $cH86 $xlRight $medium
Set-Edge $cH86 $xlTop $medium
$cH86.HorizontalAlignment = $xlCenter
$cH86.WrapText = $true

$cI86 = $ws.Range("I86")
$cI86.Value = "Saliva"
$cI86.Font.Bold = $true
Set-Edge $cI86 $xlTop $medium
$cI86.HorizontalAlignment = $xlCenter
$cI86.WrapText = $true

$cJ86 = $ws.Range("J86")
$cJ86.Value = "Henry' constant @ 310.15 K [conc_water/conc_air]"
$cJ86.Font.Bold = $true
Set-Edge $cJ86 $xlTop $medium
$cJ86.HorizontalAlignment = $xlCenter
$cJ86.WrapText = $true
$cJ86.NumberFormat = "0.00E+00"

$cK86 = $ws.Range("K86")
$cK86.Value = "Blood:Air partition coefficient [conc_blood/conc_air]"
$cK86.Font.Bold = $true
Set-Edge $cK86 $xlLeft $medium
Set-Edge $cK86 $xlRight $medium
Set-Edge $cK86 $xlTop $medium
$cK86.HorizontalAlignment = $xlCenter
$cK86.WrapText = $true
$cK86.NumberFormat = "0.00E+00"

$cL86 = $ws.Range("L86")
$cL86.Value = "Fat:Air partition coefficient  [conc_fat/conc_air]"
$cL86.Font.Bold = $true
Set-Edge $cL86 $xlRight $medium
Set-Edge $cL86 $xlTop $medium
$cL86.HorizontalAlignment = $xlCenter
$cL86.WrapText = $true
$cL86.NumberFormat = "0.00E+00"

$cM86 = $ws.Range("M86")
$cM86.Value = "Fat_Blood partition coefficient  [conc_fat/conc_blood]"
$cM86.Font.Bold = $true
Set-Edge $cM86 $xlRight $medium
Set-Edge $cM86 $xlTop $medium
$cM86.HorizontalAlignment = $xlCenter
$cM86.WrapText = $true
$cM86.NumberFormat = "0.00E+00"

# ---------------------------------------------------------------
# Update the sheet view: scrolled down to the new table, with the
# last written cell selected (mirrors the authored file).
# ---------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K87").Select()

Write-Host "done"
